# Word COM-interop script: apply the "alternate branch" edit.
#
# Target result for the single paragraph in the document:
#   Run 1 (plain):           "Hy I am creathing this file for "
#   Run 2 (Arial/11.5pt/shd): "alternate"
#   Run 3 (Arial/11.5pt/shd): " branch"

$d = $word.ActiveDocument

# Step 1: rewrite the paragraph's text and drop in two unique
# placeholder tokens marking where the newly (re)formatted runs go.
# This is a plain-text Find/Replace, so it stays inside the existing
# (unformatted) run - no run split happens here.
$d.Content.Find.Execute( `
    "Hi I am creating this file in main branch ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Hy I am creathing this file for @@ALT_RUN@@@@BR_RUN@@", 2)

# OOXML fragments (FlatOpc "single file package" form) used to replace
# the placeholder tokens. Each carries its own <w:rPr> so the inserted
# text becomes its own run with the desired direct formatting.
$altRunXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:rPr><w:t>alternate</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$brRunXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:rPr><w:t xml:space="preserve"> branch</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# InsertXML() on a non-empty Range REPLACES that range's contents, and
# the replacement text always lands at the end of its paragraph - not
# literally "in place". So to get the final reading order
# "... for " + "alternate" + " branch", resolve/replace the *second*
# placeholder ("@@BR_RUN@@") first, then the *first* one
# ("@@ALT_RUN@@"), so that the run inserted last is the one that ends
# up last in the paragraph.

$brRange = $d.Content
$brRange.Find.Execute("@@BR_RUN@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$brRange.InsertXML($brRunXml)

$altRange = $d.Content
$altRange.Find.Execute("@@ALT_RUN@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$altRange.InsertXML($altRunXml)

Write-Output "applied alternate-branch edit"
